$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D column retains text formatting (values like "30.864.16" or "1.000"
# must remain text, not be auto-converted to numbers/dates by Excel)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '30.864.16'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.890.50'
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '240.31'
$ws.Range("E5").Value = '  +2.14%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.4803'
$ws.Range("E7").Value = '  +2.23%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.2948'
$ws.Range("E8").Value = '  +6.50%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.06629'
$ws.Range("E9").Value = '  +3.79%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '18.88'
$ws.Range("E10").Value = '  +4.69%  '

$ws.Range("B11").Value = 'Litecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D11").Value = '101.07'
$ws.Range("E11").Value = '  +18.38%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07644'
$ws.Range("E12").Value = '  +2.63%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.872.65'
$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.155'
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.6595'
$ws.Range("E15").Value = '  +3.23%  '

$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").Value = '305.76'
$ws.Range("E16").Value = '  +26.64%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '30.837.58'
$ws.Range("E17").Value = '  +0.74%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '13.22'
$ws.Range("E18").Value = '  +2.74%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007640'
$ws.Range("E19").Value = '  +3.19%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.128.78'
$ws.Range("E21").Value = '  +1.34%  '

$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.193'
$ws.Range("E23").Value = '  +3.87%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '6.220'
$ws.Range("E24").Value = '  +2.35%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.369'
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '169.24'
$ws.Range("E26").Value = '  +1.97%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '20.79'
$ws.Range("E27").Value = '  +13.99%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.972'
$ws.Range("E28").Value = '  +4.13%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.1136'
$ws.Range("E29").Value = '  +11.01%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.348'
$ws.Range("E30").Value = '  -2.32%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.210'
$ws.Range("E31").Value = '  +2.72%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.010'
$ws.Range("E32").Value = '  +3.43%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.05092'
$ws.Range("E33").Value = '  +3.15%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.7519'
$ws.Range("E34").Value = '  +5.49%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.165'
$ws.Range("E35").Value = '  +0.88%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.726'
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02002'
$ws.Range("E37").Value = '  +4.70%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.701'
$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.063'
$ws.Range("E39").Value = '  +3.31%  '

$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").Value = '109.35'
$ws.Range("E40").Value = '  +3.17%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.8795'
$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4231'
$ws.Range("E43").Value = '  +2.71%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '68.88'
$ws.Range("E44").Value = '  +10.05%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '5.673'
$ws.Range("E45").Value = '  +2.08%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.389'
$ws.Range("E46").Value = '  -0.44%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1235'
$ws.Range("E47").Value = '  -0.07%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.119'
$ws.Range("E48").Value = '  +4.77%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '35.05'
$ws.Range("E49").Value = '  +4.04%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.405'
$ws.Range("E50").Value = '  +1.87%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05652'
$ws.Range("E51").Value = '  +1.39%  '
